$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# Update the "Criterio de entrada" (column C) text for rows 2-8 with the
# new entry-criteria wording for the individual plans.
$ws.Range("C2").Value = "El equipo ha completado un ciclo preeviamente."
$ws.Range("C3").Value = "El equipo ha completado un ciclo preeviamente. Cada miembro del equipo ha leído el capítulo correspondiente a su rol."
$ws.Range("C4").Value = "Se ha creado la versión final del documento de requerimientos."
$ws.Range("C5").Value = "Se ha creado el reporte de inspección del documento de requerimientos."
$ws.Range("C6").Value = "Se ha creado la versión final del documento de arquitectura."
$ws.Range("C7").Value = "Se ha implementado la funcionalidad #2 del producto."
$ws.Range("C8").Value = "El equipo ha completado los productos especificados. El equipo ha acumulado toda la información y ha completado todas las formas requeridas."

# Update the sheet view: move the active selection from H4 to C4 (this also
# clears the previous scrolled-to-D1 top-left cell override).
$ws.Range("C4").Select()
